# Updated function lists referencing lddbd_core.php and lddbd_ajax.php
# to reflect correct line locations.

$wb = $excel.ActiveWorkbook

$wsAjax = $wb.Worksheets.Item("lddbd_ajax.php")
$wsCore = $wb.Worksheets.Item("lddbd_core.php")

# --- lddbd_ajax.php: correct the LINE column (A) values ---------------
$ajaxLines = @{
    2  = 4
    3  = 13
    4  = 16
    5  = 30
    6  = 45
    7  = 56
    8  = 140
    9  = 238
    10 = 293
    11 = 379
    12 = 428
    13 = 502
    14 = 517
    15 = 550
    16 = 560
    17 = 835
    18 = 845
    19 = 861
}
foreach ($row in $ajaxLines.Keys) {
    $wsAjax.Cells.Item($row, 1).Value = $ajaxLines[$row]
}

# --- lddbd_core.php: remove the obsolete $lddbd_state_dropdown row ----
# (this entry no longer exists in the source file; deleting the row
# shifts every following entry up by one and drops the now-unused
# shared strings automatically on save)
$wsCore.Rows("6:6").Delete()

# --- lddbd_core.php: correct the LINE column (A) values for the rows
#     that shifted up -----------------------------------------------
$coreLines = @{
    6  = 36
    7  = 42
    8  = 73
    9  = 84
    10 = 94
    12 = 102
}
foreach ($row in $coreLines.Keys) {
    $wsCore.Cells.Item($row, 1).Value = $coreLines[$row]
}

# Row 11's LINE cell holds a text range reference rather than a single
# number; update it to match the corrected location.
$wsCore.Range("A11").Value = "95 to 97"

# --- restore view/selection state -------------------------------------
# Move the lddbd_core.php selection onto the new last row, then
# reselect lddbd_ajax.php (the originally active tab) on its last row
# so tab/selection state matches the corrected layout.
$wsCore.Range("A12").Select()
$wsAjax.Select()
$wsAjax.Range("A19").Select()
